$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 1491.8889  # H12
$ws.Cells.Item(12, 9).Value = 573.8333  # I12
$ws.Cells.Item(12, 10).Value = 3328  # J12
$ws.Cells.Item(12, 11).Value = 573.8333  # K12
$ws.Cells.Item(12, 12).Value = 3328  # L12
$ws.Cells.Item(12, 13).Value = -403.8333  # M12
$ws.Cells.Item(12, 14).Value = -3668  # N12

$ws.Cells.Item(40, 8).Value = 3999.913  # H40
$ws.Cells.Item(40, 10).Value = 4000  # J40
$ws.Cells.Item(40, 12).Value = 4000  # L40
$ws.Cells.Item(40, 14).Value = -4350  # N40

$ws.Cells.Item(53, 8).Value = 2001.6  # H53
$ws.Cells.Item(53, 9).Value = 903.5  # I53
$ws.Cells.Item(53, 10).Value = 3648.75  # J53
$ws.Cells.Item(53, 11).Value = 903.5  # K53
$ws.Cells.Item(53, 12).Value = 3648.75  # L53
$ws.Cells.Item(53, 13).Value = -266.5  # M53
$ws.Cells.Item(53, 14).Value = -4922.75  # N53

$ws.Cells.Item(92, 8).Value = 15625504  # H92
$ws.Cells.Item(92, 9).Value = 16667154  # I92
$ws.Cells.Item(92, 11).Value = 16667154  # K92
$ws.Cells.Item(92, 13).Value = -16665906  # M92

$ws.Cells.Item(100, 8).Value = 2803.9473  # H100
$ws.Cells.Item(100, 9).Value = 1790.3846  # I100
$ws.Cells.Item(100, 11).Value = 1790.3846  # K100
$ws.Cells.Item(100, 13).Value = -1249.3846  # M100

$ws.Cells.Item(116, 8).Value = 14198.143  # H116
$ws.Cells.Item(116, 9).Value = 16934.455  # I116
$ws.Cells.Item(116, 11).Value = 16934.455  # K116
$ws.Cells.Item(116, 13).Value = -13492.455  # M116

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20409266  # H32
$ws.Cells.Item(32, 9).Value = 21277468  # I32
$ws.Cells.Item(32, 11).Value = 21277468  # K32
$ws.Cells.Item(32, 13).Value = -21277181  # M32

$ws.Cells.Item(63, 8).Value = 4312.375  # H63
$ws.Cells.Item(63, 10).Value = 5400  # J63
$ws.Cells.Item(63, 12).Value = 5400  # L63
$ws.Cells.Item(63, 14).Value = -6772  # N63

$ws.Cells.Item(66, 8).Value = 4312.375  # H66
$ws.Cells.Item(66, 10).Value = 5400  # J66
$ws.Cells.Item(66, 12).Value = 27000  # L66
$ws.Cells.Item(66, 14).Value = -33864  # N66

$ws.Cells.Item(102, 8).Value = 1745.8334  # H102
$ws.Cells.Item(102, 9).Value = 1521.6666  # I102
$ws.Cells.Item(102, 11).Value = 1521.6666  # K102
$ws.Cells.Item(102, 13).Value = 100.3334  # M102

$ws.Cells.Item(132, 8).Value = 3876.2258  # H132
$ws.Cells.Item(132, 9).Value = 3755.8572  # I132
$ws.Cells.Item(132, 11).Value = 11267.5716  # K132
$ws.Cells.Item(132, 13).Value = -8737.571599999999  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(31, 8).Value = 0  # H31
$ws.Cells.Item(31, 10).Value = 0  # J31
$ws.Cells.Item(31, 12).ClearContents()  # L31
$ws.Cells.Item(31, 14).Value = 0  # N31

$ws.Cells.Item(94, 8).Value = 396.94116  # H94
$ws.Cells.Item(94, 9).Value = 327.14285  # I94
$ws.Cells.Item(94, 10).Value = 722.6667  # J94
$ws.Cells.Item(94, 11).Value = 327.14285  # K94
$ws.Cells.Item(94, 12).Value = 722.6667  # L94
$ws.Cells.Item(94, 13).Value = 123.85715  # M94
$ws.Cells.Item(94, 14).Value = -1624.6667  # N94

$ws.Cells.Item(99, 8).Value = 2537.375  # H99
$ws.Cells.Item(99, 9).Value = 1202.25  # I99
$ws.Cells.Item(99, 11).Value = 1202.25  # K99
$ws.Cells.Item(99, 13).Value = 295.75  # M99

$ws.Cells.Item(105, 8).Value = 2168.889  # H105
$ws.Cells.Item(105, 9).Value = 1919.3846  # I105
$ws.Cells.Item(105, 10).Value = 2817.6  # J105
$ws.Cells.Item(105, 11).Value = 1919.3846  # K105
$ws.Cells.Item(105, 12).Value = 2817.6  # L105
$ws.Cells.Item(105, 13).Value = -172.3846000000001  # M105
$ws.Cells.Item(105, 14).Value = -6311.6  # N105

$ws.Cells.Item(133, 8).Value = 97472.75  # H133
$ws.Cells.Item(133, 10).Value = 97472.75  # J133
$ws.Cells.Item(133, 12).Value = 97472.75  # L133
$ws.Cells.Item(133, 14).Value = -107592.75  # N133

$ws.Cells.Item(134, 8).Value = 4578.174  # H134
$ws.Cells.Item(134, 9).Value = 4307  # I134
$ws.Cells.Item(134, 11).Value = 12921  # K134
$ws.Cells.Item(134, 13).Value = -10386  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 28.23077  # H7
$ws.Cells.Item(7, 10).Value = 19.5  # J7
$ws.Cells.Item(7, 12).Value = 19.5  # L7
$ws.Cells.Item(7, 14).Value = -245.5  # N7

$ws.Cells.Item(16, 8).Value = 1040.3529  # H16
$ws.Cells.Item(16, 9).Value = 1121.0714  # I16
$ws.Cells.Item(16, 11).Value = 1121.0714  # K16
$ws.Cells.Item(16, 13).Value = -834.0714  # M16

$ws.Cells.Item(92, 8).Value = 59955.332  # H92
$ws.Cells.Item(92, 10).Value = 59955.332  # J92
$ws.Cells.Item(92, 12).Value = 59955.332  # L92
$ws.Cells.Item(92, 14).Value = -64947.332  # N92

$ws.Cells.Item(105, 8).Value = 1966.5714  # H105
$ws.Cells.Item(105, 9).Value = 1548.3636  # I105
$ws.Cells.Item(105, 11).Value = 1548.3636  # K105
$ws.Cells.Item(105, 13).Value = 198.6364000000001  # M105

$ws.Cells.Item(113, 8).Value = 1040.3529  # H113
$ws.Cells.Item(113, 9).Value = 1121.0714  # I113
$ws.Cells.Item(113, 11).Value = 1121.0714  # K113
$ws.Cells.Item(113, 13).Value = 1048.9286  # M113

$ws.Cells.Item(134, 8).Value = 2588.6667  # H134
$ws.Cells.Item(134, 9).Value = 2249.6667  # I134
$ws.Cells.Item(134, 11).Value = 6749.000100000001  # K134
$ws.Cells.Item(134, 13).Value = -4214.000100000001  # M134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 0  # H13
$ws.Cells.Item(13, 9).Value = 0  # I13
$ws.Cells.Item(13, 10).Value = 0  # J13
$ws.Cells.Item(13, 11).Value = 0  # K13
$ws.Cells.Item(13, 12).ClearContents()  # L13
$ws.Cells.Item(13, 13).ClearContents()  # M13
$ws.Cells.Item(13, 14).Value = 0  # N13

$ws.Cells.Item(129, 8).Value = 1558.3529  # H129
$ws.Cells.Item(129, 9).Value = 499.42856  # I129
$ws.Cells.Item(129, 10).Value = 2299.6  # J129
$ws.Cells.Item(129, 11).Value = 1498.28568  # K129
$ws.Cells.Item(129, 12).Value = 6898.799999999999  # L129
$ws.Cells.Item(129, 13).Value = 3501.71432  # M129
$ws.Cells.Item(129, 14).Value = -16898.8  # N129

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 14010.129  # H70
$ws.Cells.Item(70, 10).Value = 15409.814  # J70
$ws.Cells.Item(70, 12).Value = 15409.814  # L70
$ws.Cells.Item(70, 14).Value = -15949.814  # N70

$ws.Cells.Item(73, 8).Value = 14010.129  # H73
$ws.Cells.Item(73, 10).Value = 15409.814  # J73
$ws.Cells.Item(73, 12).Value = 15409.814  # L73
$ws.Cells.Item(73, 14).Value = -17281.814  # N73

$ws.Cells.Item(80, 8).Value = 2250  # H80
$ws.Cells.Item(80, 9).Value = 2000  # I80
$ws.Cells.Item(80, 10).Value = 3000  # J80
$ws.Cells.Item(80, 11).Value = 2000  # K80
$ws.Cells.Item(80, 12).Value = 3000  # L80
$ws.Cells.Item(80, 13).Value = -1002  # M80
$ws.Cells.Item(80, 14).Value = -4996  # N80

$ws.Cells.Item(83, 8).Value = 2250  # H83
$ws.Cells.Item(83, 9).Value = 2000  # I83
$ws.Cells.Item(83, 10).Value = 3000  # J83
$ws.Cells.Item(83, 11).Value = 10000  # K83
$ws.Cells.Item(83, 12).Value = 15000  # L83
$ws.Cells.Item(83, 13).Value = -5008  # M83
$ws.Cells.Item(83, 14).Value = -24984  # N83

$ws.Cells.Item(122, 8).Value = 2989.375  # H122
$ws.Cells.Item(122, 9).Value = 3433.6  # I122
$ws.Cells.Item(122, 10).Value = 2249  # J122
$ws.Cells.Item(122, 11).Value = 10300.8  # K122
$ws.Cells.Item(122, 12).Value = 6747  # L122
$ws.Cells.Item(122, 13).Value = -7850.799999999999  # M122
$ws.Cells.Item(122, 14).Value = -11647  # N122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 1408  # H40
$ws.Cells.Item(40, 9).Value = 1408  # I40
$ws.Cells.Item(40, 10).Value = 0  # J40
$ws.Cells.Item(40, 11).Value = 1408  # K40
$ws.Cells.Item(40, 12).Value = 0  # L40
$ws.Cells.Item(40, 13).ClearContents()  # M40
$ws.Cells.Item(40, 14).Value = -1272  # N40

$ws.Cells.Item(46, 8).Value = 3163.3142  # H46
$ws.Cells.Item(46, 9).Value = 1664  # I46
$ws.Cells.Item(46, 10).Value = 3538.1428  # J46
$ws.Cells.Item(46, 11).Value = 1664  # K46
$ws.Cells.Item(46, 12).Value = 3538.1428  # L46
$ws.Cells.Item(46, 13).Value = -1476  # M46
$ws.Cells.Item(46, 14).Value = -3914.1428  # N46

$ws.Cells.Item(93, 8).Value = 2322.9167  # H93
$ws.Cells.Item(93, 9).Value = 2229.1538  # I93
$ws.Cells.Item(93, 10).Value = 2433.7273  # J93
$ws.Cells.Item(93, 11).Value = 2229.1538  # K93
$ws.Cells.Item(93, 12).Value = 2433.7273  # L93
$ws.Cells.Item(93, 13).Value = -981.1538  # M93
$ws.Cells.Item(93, 14).Value = -4929.7273  # N93

$ws.Cells.Item(100, 8).Value = 1355.5555  # H100
$ws.Cells.Item(100, 9).Value = 1166.6666  # I100
$ws.Cells.Item(100, 10).Value = 1733.3334  # J100
$ws.Cells.Item(100, 11).Value = 1166.6666  # K100
$ws.Cells.Item(100, 12).Value = 1733.3334  # L100
$ws.Cells.Item(100, 13).Value = -625.6666  # M100
$ws.Cells.Item(100, 14).Value = -2815.3334  # N100

$ws.Cells.Item(122, 8).Value = 12552.857  # H122
$ws.Cells.Item(122, 9).Value = 12362.8125  # I122
$ws.Cells.Item(122, 10).Value = 13161  # J122
$ws.Cells.Item(122, 11).Value = 37088.4375  # K122
$ws.Cells.Item(122, 12).Value = 39483  # L122
$ws.Cells.Item(122, 13).Value = -34638.4375  # M122
$ws.Cells.Item(122, 14).Value = -44383  # N122

$ws.Cells.Item(125, 8).Value = 78990  # H125
$ws.Cells.Item(125, 10).Value = 78990  # J125
$ws.Cells.Item(125, 12).Value = 78990  # L125
$ws.Cells.Item(125, 14).Value = -88830  # N125

$ws.Cells.Item(127, 8).Value = 21277.777  # H127
$ws.Cells.Item(127, 10).Value = 21277.777  # J127
$ws.Cells.Item(127, 12).Value = 21277.777  # L127
$ws.Cells.Item(127, 14).Value = -31197.777  # N127

$ws.Cells.Item(130, 8).Value = 0  # H130
$ws.Cells.Item(130, 10).Value = 0  # J130
$ws.Cells.Item(130, 12).ClearContents()  # L130
$ws.Cells.Item(130, 14).Value = 0  # N130

$ws.Cells.Item(132, 8).Value = 1325  # H132
$ws.Cells.Item(132, 9).Value = 1325  # I132
$ws.Cells.Item(132, 11).Value = 3975  # K132
$ws.Cells.Item(132, 13).Value = -1445  # M132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4649.923  # H62
$ws.Cells.Item(62, 9).Value = 3124.75  # I62
$ws.Cells.Item(62, 10).Value = 5327.778  # J62
$ws.Cells.Item(62, 11).Value = 3124.75  # K62
$ws.Cells.Item(62, 12).Value = 5327.778  # L62
$ws.Cells.Item(62, 13).Value = -2500.75  # M62
$ws.Cells.Item(62, 14).Value = -6575.778  # N62

$ws.Cells.Item(65, 8).Value = 4649.923  # H65
$ws.Cells.Item(65, 9).Value = 3124.75  # I65
$ws.Cells.Item(65, 10).Value = 5327.778  # J65
$ws.Cells.Item(65, 11).Value = 15623.75  # K65
$ws.Cells.Item(65, 12).Value = 26638.89  # L65
$ws.Cells.Item(65, 13).Value = -12503.75  # M65
$ws.Cells.Item(65, 14).Value = -32878.89  # N65

$ws.Cells.Item(107, 8).Value = 365.09525  # H107
$ws.Cells.Item(107, 9).Value = 268.82352  # I107
$ws.Cells.Item(107, 11).Value = 806.47056  # K107
$ws.Cells.Item(107, 13).Value = 1113.52944  # M107

$ws.Cells.Item(132, 8).Value = 2943.0557  # H132
$ws.Cells.Item(132, 9).Value = 2838.3333  # I132
$ws.Cells.Item(132, 11).Value = 8514.999899999999  # K132
$ws.Cells.Item(132, 13).Value = -5984.999899999999  # M132

$ws.Cells.Item(136, 8).Value = 2889.4707  # H136
$ws.Cells.Item(136, 9).Value = 1625.1538  # I136
$ws.Cells.Item(136, 11).Value = 4875.4614  # K136
$ws.Cells.Item(136, 13).Value = -2325.4614  # M136
